$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")
Write-Host ("Row count: " + $lo.ListRows.Count)
$lo.ListRows.Add(151) | Out-Null
Write-Host ("After rows: " + $lo.ListRows.Count)
Write-Host ("After range: " + $lo.Range.Address())
